# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.395.68"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "1.868.48"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'339.12"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.4699"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").Value = "'0.3957"
$ws.Range("E8").Value = "  +3.77%  "
$ws.Range("D9").Value = "'47.33"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").Value = "'0.08008"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").Value = "'0.9996"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("D12").Value = "'21.84"
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").Value = "1.862.49"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'5.987"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "'7.225"
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "'91.24"
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "'0.00001042"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'0.06630"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "'17.56"
$ws.Range("E20").Value = "  +3.79%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "28.416.38"
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("D23").Value = "'5.451"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").Value = "'11.03"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").Value = "2.098.45"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "'160.45"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").Value = "'19.76"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "'2.129"
$ws.Range("E29").Value = "  +3.08%  "
$ws.Range("D30").Value = "'5.502"
$ws.Range("E30").Value = "  +3.79%  "
$ws.Range("D31").Value = "'119.97"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "'0.9653"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").Value = "'0.09480"
$ws.Range("D34").Value = "'3.570"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.347"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.372"
$ws.Range("E36").Value = "  +4.11%  "
$ws.Range("D37").Value = "'0.06083"
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("D38").Value = "'0.02244"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "'8.372"
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("D40").Value = "'1.185"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").Value = "'0.5937"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Value = "'0.1867"
$ws.Range("E43").Value = "  +1.63%  "
$ws.Range("D44").Value = "'10.32"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").Value = "'1.291"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").Value = "'0.5576"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "'12.13"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Value = "'1.953"
$ws.Range("E48").Value = "  +4.84%  "
$ws.Range("D49").Value = "'0.06854"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").Value = "'2.064"
$ws.Range("E50").Value = "  +16.65%  "
$ws.Range("D51").Value = "'111.32"
$ws.Range("E51").Value = "  +1.51%  "
